# Update Suryakumar Yadav's per-innings batting log (runs/balls/fours/sixes)
# to reflect the latest activity pulled in from the Excel form. Values are
# stored as text in the sheet, so each assignment is prefixed with an
# apostrophe to force a text literal (avoiding Excel's automatic numeric
# coercion when a cell is set to a numeric-looking string).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'0"
$ws.Range("D2").Value = "'4"
$ws.Range("E2").Value = "'0"

$ws.Range("C3").Value = "'12"
$ws.Range("D3").Value = "'11"
$ws.Range("E3").Value = "'1"
$ws.Range("F3").Value = "'0"

$ws.Range("D4").Value = "'10"
$ws.Range("E4").Value = "'1"

$ws.Range("C5").Value = "'36"
$ws.Range("D5").Value = "'29"
$ws.Range("E5").Value = "'5"

$ws.Range("C6").Value = "'79"
$ws.Range("D6").Value = "'43"
$ws.Range("E6").Value = "'10"
$ws.Range("F6").Value = "'3"

$ws.Range("C7").Value = "'40"
$ws.Range("D7").Value = "'26"
$ws.Range("E7").Value = "'4"
$ws.Range("F7").Value = "'1"

$ws.Range("C8").Value = "'51"
$ws.Range("D8").Value = "'38"
$ws.Range("E8").Value = "'6"
$ws.Range("F8").Value = "'2"

$ws.Range("C9").Value = "'79"
$ws.Range("D9").Value = "'47"
$ws.Range("E9").Value = "'11"
$ws.Range("F9").Value = "'2"

$ws.Range("C11").Value = "'17"
$ws.Range("D11").Value = "'16"
$ws.Range("E11").Value = "'2"
$ws.Range("F11").Value = "'0"

$ws.Range("C13").Value = "'47"
$ws.Range("D13").Value = "'28"
$ws.Range("F13").Value = "'1"

$ws.Range("C14").Value = "'27"
$ws.Range("D14").Value = "'18"
$ws.Range("F14").Value = "'0"

$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'2"
$ws.Range("E15").Value = "'0"

$ws.Range("C16").Value = "'10"
$ws.Range("D16").Value = "'7"
$ws.Range("E16").Value = "'2"
$ws.Range("F16").Value = "'0"
